$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.010597825050354
$ws.Range("B1").Value = 2.124265193939209
$ws.Range("C1").Value = 6.063383102416992
$ws.Range("D1").Value = 1.280569672584534
$ws.Range("E1").Value = 1.265876412391663
